# Adapt column header formatting to respective input file names:
#   "<header>_old" -> "<header>_FV2310"
#   "<header>_new" -> "<header>_FV2404"
# Also turn the data range into an Excel Table ("Table1") and freeze the
# header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fv2310Headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

# Columns A..J hold the "_old" -> "_FV2310" headers, column K holds "diff",
# columns L..U hold the "_new" -> "_FV2404" headers.
for ($i = 0; $i -lt $fv2310Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2310Headers[$i]
}

for ($i = 0; $i -lt $fv2404Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2404Headers[$i]
}

# Turn the data range into an Excel Table (ListObject) named "Table1".
$tableRange = $ws.Range("A1:U89")
$listObject = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$listObject.Name = "Table1"

# Freeze the header row (pane split after row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
